# Generate Report for Archive
#
# Refresh the localization status report: the two files that were still
# "Ready for handoff" but have since entered translation
# (3e7a8619-6e35-4aee-8675-9767fe6d503c.md and
#  94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.md) now show status "In Translation",
# and the rows are re-sorted by Status (Handed back < In Translation <
# Ready for handoff), which swaps the 94abd0c5 / 8692af51 rows on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 3 (3e7a8619...md): status flips to "In Translation" in both locale
# columns; date/file columns are unchanged.
$ws.Range("E3").Value2 = "In Translation"
$ws.Range("F3").Value2 = "In Translation"

# Row 4 becomes the 94abd0c5 file (was the 8692af51 file), now "In
# Translation", dated 2016-09-04 00:45:25.
$ws.Range("A4").Value2 = "94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.md"
$ws.Range("B4").Value2 = "e2e\94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.md"
$ws.Range("E4").Value2 = "In Translation"
$ws.Range("F4").Value2 = "In Translation"
$ws.Range("G4").Value2 = "2016-09-04 00:45:25"

# Row 5 becomes the 8692af51 file (was the 94abd0c5 file), still "Ready
# for handoff", dated 2016-09-04 00:44:20.
$ws.Range("A5").Value2 = "8692af51-04ec-4baa-951a-ed5ad53d6d4c.md"
$ws.Range("B5").Value2 = "e2e\8692af51-04ec-4baa-951a-ed5ad53d6d4c.md"
$ws.Range("G5").Value2 = "2016-09-04 00:44:20"

# The hyperlinks on column B keep pointing at their original commit URLs,
# but their displayed text must track the cell text above.
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$4') {
        $h.TextToDisplay = "e2e\94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.md"
    } elseif ($addr -eq '$B$5') {
        $h.TextToDisplay = "e2e\8692af51-04ec-4baa-951a-ed5ad53d6d4c.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C3").Value2 = "In Translation"

$ws.Range("A4").Value2 = "94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.md"
$ws.Range("C4").Value2 = "In Translation"
$ws.Range("G4").Value2 = "94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.7688bba6240a108f2fa5d856f7eb7d86080df54c.zh-cn.xlf"
$ws.Range("H4").Value2 = "2016-09-04 00:45:21"

$ws.Range("A5").Value2 = "8692af51-04ec-4baa-951a-ed5ad53d6d4c.md"
$ws.Range("G5").Value2 = "8692af51-04ec-4baa-951a-ed5ad53d6d4c.dbf421b1cd1fa08ef5c60aeea814025a0ee0e740.zh-cn.xlf"
$ws.Range("H5").Value2 = "2016-09-04 00:44:14"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$4') {
        $h.TextToDisplay = "94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.md"
    } elseif ($addr -eq '$A$5') {
        $h.TextToDisplay = "8692af51-04ec-4baa-951a-ed5ad53d6d4c.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C3").Value2 = "In Translation"

$ws.Range("A4").Value2 = "94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.md"
$ws.Range("C4").Value2 = "In Translation"
$ws.Range("G4").Value2 = "94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.7688bba6240a108f2fa5d856f7eb7d86080df54c.de-de.xlf"
$ws.Range("H4").Value2 = "2016-09-04 00:45:25"

$ws.Range("A5").Value2 = "8692af51-04ec-4baa-951a-ed5ad53d6d4c.md"
$ws.Range("G5").Value2 = "8692af51-04ec-4baa-951a-ed5ad53d6d4c.dbf421b1cd1fa08ef5c60aeea814025a0ee0e740.de-de.xlf"
$ws.Range("H5").Value2 = "2016-09-04 00:44:20"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$4') {
        $h.TextToDisplay = "94abd0c5-ae9b-4cd1-a9ec-2b68c61e9845.md"
    } elseif ($addr -eq '$A$5') {
        $h.TextToDisplay = "8692af51-04ec-4baa-951a-ed5ad53d6d4c.md"
    }
}
